$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A130").Value = "hassan"
$ws.Range("B130").Value = "hasan"
$ws.Range("A131").Value = "hussein"
$ws.Range("B131").Value = "hussen, husein, hussain, husain"

$ws.Range("B132").Select()
